$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data set: jump/break point handling order changed,
# inserting earlier observations and re-deriving later ones, shifting all rows.
$dates = @(
    "20150914",
    "20150918",
    "20151125",
    "20151211",
    "20151223",
    "20160112",
    "20160222",
    "20160229",
    "20160321",
    "20160329",
    "20160412",
    "20160520",
    "20160617",
    "20160624",
    "20160727",
    "20160801",
    "20161012",
    "20161021",
    "20161118",
    "20170116",
    "20170208",
    "20170511",
    "20170710",
    "20170724",
    "20170907",
    "20170922",
    "20171011",
    "20171103",
    "20171117",
    "20171206",
    "20171226",
    "20180207",
    "20180314",
    "20180326",
    "20180402",
    "20180423",
    "20180511",
    "20180709",
    "20180713",
    "20180903",
    "20180926",
    "20181019",
    "20181204",
    "20181218",
    "20190325",
    "20190506",
    "20190528",
    "20190610",
    "20190621",
    "20190722",
    "20190731",
    "20190806",
    "20190910",
    "20191021",
    "20191105",
    "20191202",
    "20200327",
    "20200402",
    "20200416",
    "20200428",
    "20200617",
    "20200624"
)

$values = @(
    821.0,
    780.0,
    1254.0,
    1055.0,
    1283.0,
    839.0,
    965.0,
    745.0,
    907.0,
    837.0,
    1029.0,
    787.0,
    941.0,
    849.0,
    985.0,
    862.0,
    1054.0,
    971.0,
    1151.0,
    786.0,
    910.0,
    643.0,
    711.0,
    612.0,
    787.0,
    695.0,
    756.0,
    624.0,
    678.0,
    566.0,
    837.0,
    528.0,
    746.0,
    594.0,
    791.0,
    609.0,
    710.0,
    477.0,
    528.0,
    415.0,
    467.0,
    323.0,
    676.0,
    470.0,
    847.0,
    544.0,
    657.0,
    566.0,
    635.0,
    532.0,
    566.0,
    468.0,
    583.0,
    478.0,
    538.0,
    468.0,
    666.0,
    526.0,
    625.0,
    523.0,
    620.0,
    570.0
)

$rowCount = $dates.Length

# Force column A to text so the numeric-looking dates are not coerced to numbers.
$dateRange = $ws.Range("A1:A$rowCount")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Drop the temporary text format so the cell style matches the original (default) style.
$dateRange.ClearFormats()
